$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 479, pushing the old rows 479:543 down to 480:544.
$ws.Rows.Item(479).Insert()

# Copy the row right below (which now holds what used to be row 479) up into the
# freshly-inserted row 479 so every column that doesn't change keeps identical
# content/format.
$ws.Rows.Item(480).Copy()
$ws.Rows.Item(479).PasteSpecial()

# Now overwrite just the columns that actually carry new data for this new record.
$ws.Range("D479").Value = 45131
$ws.Range("J479").Value = 125
$ws.Range("K479").Value = 8000
$ws.Range("L479").Value = 8000
$ws.Range("M479").Value = 8000
$ws.Range("P479").Value = 667
